# algo_make_setup.xlsx - "Adding folder for make_setup.py"
#
# The sheet is a flattened directory listing (parent / folders / type / path).
# It grows from 9 data rows (one row per folder *and* a second combined
# row per text file) to 15 data rows (one row per item, with an explicit
# "type" column of FOLDER/TEXT), plus a header row. Below we rewrite the
# whole data block and restyle it to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Header row - only D1 actually changes text (files -> type)
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "parent"
$ws.Range("C1").Value = "folders"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "path"
$ws.Range("G1").Value = "path"

# ---------------------------------------------------------------------
# 2. Data rows 2-16 (parent level, parent folder name, item name, type, path)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "folder"
$ws.Range("C2").Value = "folder1"
$ws.Range("D2").Value = "FOLDER"
$ws.Range("E2").Value = "folder"
$ws.Range("G2").Value = "C:\Users\Rohan\Desktop\folder"

$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "folder"
$ws.Range("C3").Value = "folder2"
$ws.Range("D3").Value = "FOLDER"
$ws.Range("E3").Value = "folder"

$ws.Range("A4").Value = 0
$ws.Range("B4").Value = "folder"
$ws.Range("C4").Value = "text1"
$ws.Range("D4").Value = "TEXT"
$ws.Range("E4").Value = "folder"

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = "folder"
$ws.Range("C5").Value = "text2"
$ws.Range("D5").Value = "TEXT"
$ws.Range("E5").Value = "folder"

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "folder1"
$ws.Range("C6").Value = "folder3"
$ws.Range("D6").Value = "FOLDER"
$ws.Range("E6").Value = "folder\folder1"

$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "folder1"
$ws.Range("C7").Value = "text3"
$ws.Range("D7").Value = "TEXT"
$ws.Range("E7").Value = "folder\folder1"

$ws.Range("A8").Value = 1
$ws.Range("B8").Value = "folder1"
$ws.Range("C8").Value = "text4"
$ws.Range("D8").Value = "TEXT"
$ws.Range("E8").Value = "folder\folder1"

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "folder2"
$ws.Range("C9").Value = "folder4"
$ws.Range("D9").Value = "FOLDER"
$ws.Range("E9").Value = "folder\folder2"

$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "folder2"
$ws.Range("C10").Value = "text5"
$ws.Range("D10").Value = "TEXT"
$ws.Range("E10").Value = "folder\folder2"

$ws.Range("A11").Value = 2
$ws.Range("B11").Value = "folder4"
$ws.Range("C11").Value = "folder5"
$ws.Range("D11").Value = "FOLDER"
$ws.Range("E11").Value = "folder\folder2\folder3"

$ws.Range("A12").Value = 2
$ws.Range("B12").Value = "folder4"
$ws.Range("C12").Value = "folder6"
$ws.Range("D12").Value = "FOLDER"
$ws.Range("E12").Value = "folder\folder2\folder3"

$ws.Range("A13").Value = 2
$ws.Range("B13").Value = "folder4"
$ws.Range("C13").Value = "text6"
$ws.Range("D13").Value = "TEXT"
$ws.Range("E13").Value = "folder\folder2\folder3"

$ws.Range("A14").Value = 2
$ws.Range("B14").Value = "folder4"
$ws.Range("C14").Value = "text7"
$ws.Range("D14").Value = "TEXT"
$ws.Range("E14").Value = "folder\folder2\folder3"

$ws.Range("A15").Value = 3
$ws.Range("B15").Value = "folder5"
$ws.Range("C15").Value = "text8"
$ws.Range("D15").Value = "TEXT"
$ws.Range("E15").Value = "folder\folder2\folder4\folder5"

$ws.Range("A16").Value = 3
$ws.Range("B16").Value = "folder6"
$ws.Range("E16").Value = "folder\folder2\folder4\folder6"

# ---------------------------------------------------------------------
# 3. Re-apply the existing banded-fill row styles to their new rows.
#    Each of these colours already exists in the workbook's style table,
#    so we copy the formatting from a row that still carries it rather
#    than re-creating it (keeps the style table from growing). The copy
#    order below is chosen so a row is always used as a formatting
#    "source" before it gets overwritten as a "destination".
# ---------------------------------------------------------------------
$ws.Range("A7:E7").Copy() | Out-Null
$ws.Range("A11:E14").PasteSpecial(-4122) | Out-Null

$ws.Range("A9:E9").Copy() | Out-Null
$ws.Range("A15:E15").PasteSpecial(-4122) | Out-Null

$ws.Range("A6:E6").Copy() | Out-Null
$ws.Range("A9:E10").PasteSpecial(-4122) | Out-Null

$ws.Range("A4:E4").Copy() | Out-Null
$ws.Range("A6:E8").PasteSpecial(-4122) | Out-Null

$ws.Range("A2:E2").Copy() | Out-Null
$ws.Range("A2:E5").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Row 16 is the new final row - give it a fresh blue fill while
#    reusing the workbook's existing thin box border (copied from G2,
#    which already carries that border with no fill).
# ---------------------------------------------------------------------
$ws.Range("G2").Copy() | Out-Null
$ws.Range("A16:E16").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("A16:E16").Interior.Color = 15773696

# ---------------------------------------------------------------------
# 5. Match the author's final selection in the saved file.
# ---------------------------------------------------------------------
$ws.Range("D16").Select()
